# Update "想去人数" (want-to-go count) values in column F for the
# "展览" (sheet1) and "全部类型" (sheet4) worksheets, matching the
# newly generated data output.

$wb = $excel.ActiveWorkbook

# Map: worksheet name -> { row -> new value }
$updates = @{
    "展览"     = @{ 2 = 6808; 4 = 434; 9 = 100; 12 = 20; 13 = 176; 16 = 1623; 17 = 25; 18 = 3430; 22 = 2062; 23 = 161; 28 = 9; 29 = 139 }
    "全部类型" = @{ 2 = 6808; 4 = 434; 10 = 100; 13 = 20; 14 = 176; 17 = 1623; 18 = 25; 19 = 3430; 23 = 2062; 24 = 161; 29 = 9; 30 = 139 }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Range("F$row").Value = $rows[$row]
    }
}
